$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "47.280.78"
$ws.Range("E2").Value = "  +4.98%  "
# Row 3
$ws.Range("D3").Value = "2.646.55"
$ws.Range("E3").Value = "  +11.50%  "
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.36%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.57%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "106.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +12.28%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.616"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +10.25%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.16%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.605"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +19.97%  "
# Row 10
$ws.Range("E10").Value = "  +15.55%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "56.19"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.99%  "
# Row 12
$ws.Range("E12").Value = "  +10.30%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.49"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +21.44%  "
# Row 14
$ws.Range("D14").Value = "3.050.87"
$ws.Range("E14").Value = "  +11.27%  "
# Row 15
$ws.Range("E15").Value = "  +2.75%  "
# Row 16
$ws.Range("D16").Value = "2.669.33"
$ws.Range("E16").Value = "  +12.60%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.950"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +14.91%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "15.48"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +10.64%  "
# Row 19
$ws.Range("D19").Value = "47.829.92"
$ws.Range("E19").Value = "  +6.27%  "
# Row 20
$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.66"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +9.32%  "
# Row 21
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000105"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +11.97%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.89"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +12.97%  "
# Row 23
$ws.Range("E23").Value = "  +9.72%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "275.86"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +15.28%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.15"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +13.60%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "30.95"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +47.95%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +17.85%  "
# Row 28
$ws.Range("E28").Value = "  +0.11%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.06"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.88%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.87"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +13.20%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "41.49"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +9.92%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.32"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.50%  "
# Row 33
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.36"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +16.88%  "
# Row 34
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.82"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.02%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.31"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +19.22%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0862"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +13.02%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.88"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.40%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "152.96"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.59%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.126"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +12.02%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.126"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.05%  "
# Row 41
$ws.Range("B41").Value = "Celestia"
$ws.Range("C41").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.20"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +15.37%  "
# Row 42
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "23.53"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +57.29%  "
# Row 43
$ws.Range("E43").Value = "  +15.85%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.75"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +17.57%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0334"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +13.11%  "
# Row 46
$ws.Range("D46").Value = "2.087.32"
$ws.Range("E46").Value = "  +7.72%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "99.64"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +10.78%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.998"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.01%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "115.58"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +15.55%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.88"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +10.10%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.35"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +8.72%  "
